# This script applies the data refresh for the "cryptos" price/volume table
# (GitHub Actions scheduled update), matching the supplied OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36 / Row 37: coins swapped position (NEARProtocol now ranks above ImmutableX) ---
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'4.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.04%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.30%  "

# --- Price (column D) and Volume(1h) (column E) updates for all other rows ---
$ws.Range("D2").Value = "60.847.90"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "2.412.27"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'569.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'138.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "2.394.33"
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "'25.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "'0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").Value = "2.860.38"
$ws.Range("D17").Value = "60.794.92"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "2.398.74"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "'7.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.77%  "
$ws.Range("D20").Value = "'10.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "'322.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").Value = "'6.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("D26").Value = "'64.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'577.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "'8.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.25%  "
$ws.Range("D29").Value = "2.536.79"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("D31").Value = "'7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -5.25%  "
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D38").Value = "'150.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "'0.366"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Value = "'18.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'5.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("D44").Value = "'41.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("E45").Value = "  -6.72%  "
$ws.Range("D46").Value = "'141.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "'3.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").Value = "'0.584"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("D50").Value = "'19.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("E51").Value = "  -3.44%  "
